$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 83334470
$ws.Cells.Item(18, 9).Value = 83334470
$ws.Cells.Item(18, 11).Value = 83334470
$ws.Cells.Item(18, 13).Value = -83334186
$ws.Cells.Item(33, 8).Value = 1868.5555
$ws.Cells.Item(33, 10).Value = 930.6667
$ws.Cells.Item(33, 12).Value = 930.6667
$ws.Cells.Item(33, 14).Value = -1388.6667
$ws.Cells.Item(41, 8).Value = 394.06668
$ws.Cells.Item(41, 9).Value = 217.94737
$ws.Cells.Item(41, 10).Value = 698.2727
$ws.Cells.Item(41, 11).Value = 217.94737
$ws.Cells.Item(41, 12).Value = 698.2727
$ws.Cells.Item(41, 13).Value = 222.05263
$ws.Cells.Item(41, 14).Value = -1578.2727
$ws.Cells.Item(43, 8).Value = 2270.625
$ws.Cells.Item(43, 9).Value = 2317.8572
$ws.Cells.Item(43, 10).Value = 1940
$ws.Cells.Item(43, 11).Value = 2317.8572
$ws.Cells.Item(43, 12).Value = 1940
$ws.Cells.Item(43, 13).Value = -2248.8572
$ws.Cells.Item(43, 14).Value = -2078
$ws.Cells.Item(45, 8).Value = 1191.8572
$ws.Cells.Item(45, 9).Value = 981.6667
$ws.Cells.Item(45, 10).Value = 1349.5
$ws.Cells.Item(45, 11).Value = 2945.0001
$ws.Cells.Item(45, 12).Value = 4048.5
$ws.Cells.Item(45, 13).Value = -2753.0001
$ws.Cells.Item(45, 14).Value = -4432.5
$ws.Cells.Item(64, 8).Value = 4699.8887
$ws.Cells.Item(64, 9).Value = 3574.75
$ws.Cells.Item(64, 11).Value = 3574.75
$ws.Cells.Item(64, 13).Value = -3326.75
$ws.Cells.Item(67, 8).Value = 4699.8887
$ws.Cells.Item(67, 9).Value = 3574.75
$ws.Cells.Item(67, 11).Value = 3574.75
$ws.Cells.Item(67, 13).Value = -2716.75
$ws.Cells.Item(74, 8).Value = 4821.3335
$ws.Cells.Item(76, 8).Value = 183570.72
$ws.Cells.Item(76, 9).Value = 361666.34
$ws.Cells.Item(76, 10).Value = 49999
$ws.Cells.Item(76, 11).Value = 361666.34
$ws.Cells.Item(76, 12).Value = 49999
$ws.Cells.Item(76, 13).Value = -361351.34
$ws.Cells.Item(76, 14).Value = -50629
$ws.Cells.Item(77, 8).Value = 4821.3335
$ws.Cells.Item(79, 8).Value = 183570.72
$ws.Cells.Item(79, 9).Value = 361666.34
$ws.Cells.Item(79, 10).Value = 49999
$ws.Cells.Item(79, 11).Value = 361666.34
$ws.Cells.Item(79, 12).Value = 49999
$ws.Cells.Item(79, 13).Value = -360574.34
$ws.Cells.Item(79, 14).Value = -52183
$ws.Cells.Item(86, 8).Value = 142859470
$ws.Cells.Item(86, 9).Value = 200001840
$ws.Cells.Item(86, 10).Value = 3521
$ws.Cells.Item(86, 11).Value = 200001840
$ws.Cells.Item(86, 12).Value = 3521
$ws.Cells.Item(86, 13).Value = -200000717
$ws.Cells.Item(86, 14).Value = -5767
$ws.Cells.Item(89, 8).Value = 142859470
$ws.Cells.Item(89, 9).Value = 200001840
$ws.Cells.Item(89, 10).Value = 3521
$ws.Cells.Item(89, 11).Value = 1000009200
$ws.Cells.Item(89, 12).Value = 17605
$ws.Cells.Item(89, 13).Value = -1000003584
$ws.Cells.Item(89, 14).Value = -28837
$ws.Cells.Item(92, 8).Value = 23809902
$ws.Cells.Item(92, 10).Value = 98.5
$ws.Cells.Item(92, 12).Value = 98.5
$ws.Cells.Item(92, 14).Value = -2594.5
$ws.Cells.Item(94, 8).Value = 14287444
$ws.Cells.Item(94, 9).Value = 15874937
$ws.Cells.Item(94, 11).Value = 15874937
$ws.Cells.Item(94, 13).Value = -15874486
$ws.Cells.Item(100, 8).Value = 1826.8572
$ws.Cells.Item(100, 9).Value = 1660
$ws.Cells.Item(100, 10).Value = 1879
$ws.Cells.Item(100, 11).Value = 1660
$ws.Cells.Item(100, 12).Value = 1879
$ws.Cells.Item(100, 13).Value = -1119
$ws.Cells.Item(100, 14).Value = -2961
$ws.Cells.Item(106, 8).Value = 150657.72
$ws.Cells.Item(106, 9).Value = 252924.5
$ws.Cells.Item(106, 10).Value = 14302
$ws.Cells.Item(106, 11).Value = 252924.5
$ws.Cells.Item(106, 12).Value = 14302
$ws.Cells.Item(106, 13).Value = -252293.5
$ws.Cells.Item(106, 14).Value = -15564
$ws.Cells.Item(107, 8).Value = 2267.7942
$ws.Cells.Item(107, 9).Value = 1513.5834
$ws.Cells.Item(107, 10).Value = 4077.9
$ws.Cells.Item(107, 11).Value = 1513.5834
$ws.Cells.Item(107, 12).Value = 4077.9
$ws.Cells.Item(107, 13).Value = 406.4166
$ws.Cells.Item(107, 14).Value = -7917.9
$ws.Cells.Item(109, 8).Value = 51499.668
$ws.Cells.Item(109, 10).Value = 51499.668
$ws.Cells.Item(109, 12).Value = 51499.668
$ws.Cells.Item(109, 14).Value = -54273.668
$ws.Cells.Item(116, 8).Value = 5224.2856
$ws.Cells.Item(116, 10).Value = 4291.24
$ws.Cells.Item(116, 12).Value = 4291.24
$ws.Cells.Item(116, 14).Value = -11175.24
$ws.Cells.Item(132, 8).Value = 9311.543
$ws.Cells.Item(132, 9).Value = 6266.0454
$ws.Cells.Item(132, 10).Value = 14465.462
$ws.Cells.Item(132, 11).Value = 18798.1362
$ws.Cells.Item(132, 12).Value = 43396.386
$ws.Cells.Item(132, 13).Value = -16268.1362
$ws.Cells.Item(132, 14).Value = -48456.386
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
$ws.Cells.Item(140, 8).Value = 98162
$ws.Cells.Item(140, 10).Value = 148975
$ws.Cells.Item(140, 12).Value = 148975
$ws.Cells.Item(140, 14).Value = -159335

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 698585.06
$ws.Cells.Item(2, 9).Value = 1021935.94
$ws.Cells.Item(2, 10).Value = 2137.077
$ws.Cells.Item(2, 11).Value = 1021935.94
$ws.Cells.Item(2, 12).Value = 2137.077
$ws.Cells.Item(2, 13).Value = -1021822.94
$ws.Cells.Item(2, 14).Value = -2363.077
$ws.Cells.Item(8, 8).Value = 5428.4287
$ws.Cells.Item(8, 9).Value = 2499.5
$ws.Cells.Item(8, 10).Value = 6600
$ws.Cells.Item(8, 11).Value = 2499.5
$ws.Cells.Item(8, 12).Value = 6600
$ws.Cells.Item(8, 13).Value = -2355.5
$ws.Cells.Item(8, 14).Value = -6888
$ws.Cells.Item(10, 8).Value = 4299.5
$ws.Cells.Item(10, 9).Value = 4119.4
$ws.Cells.Item(10, 11).Value = 4119.4
$ws.Cells.Item(10, 13).Value = -3949.4
$ws.Cells.Item(11, 8).Value = 2133
$ws.Cells.Item(11, 9).Value = 899.5
$ws.Cells.Item(11, 11).Value = 899.5
$ws.Cells.Item(11, 13).Value = -755.5
$ws.Cells.Item(14, 8).Value = 355.2857
$ws.Cells.Item(14, 9).Value = 334.33334
$ws.Cells.Item(14, 10).Value = 371
$ws.Cells.Item(14, 11).Value = 334.33334
$ws.Cells.Item(14, 12).Value = 371
$ws.Cells.Item(14, 13).Value = -159.33334
$ws.Cells.Item(14, 14).Value = -721
$ws.Cells.Item(16, 8).Value = 1192.1
$ws.Cells.Item(16, 9).Value = 846.7143
$ws.Cells.Item(16, 10).Value = 1998
$ws.Cells.Item(16, 11).Value = 846.7143
$ws.Cells.Item(16, 12).Value = 1998
$ws.Cells.Item(16, 13).Value = -559.7143
$ws.Cells.Item(16, 14).Value = -2572
$ws.Cells.Item(17, 8).Value = 1499.5
$ws.Cells.Item(17, 9).Value = 999
$ws.Cells.Item(17, 11).Value = 999
$ws.Cells.Item(17, 13).Value = -826
$ws.Cells.Item(18, 8).Value = 11474.667
$ws.Cells.Item(18, 9).Value = 14974.5
$ws.Cells.Item(18, 10).Value = 9724.75
$ws.Cells.Item(18, 11).Value = 14974.5
$ws.Cells.Item(18, 12).Value = 9724.75
$ws.Cells.Item(18, 13).Value = -14652.5
$ws.Cells.Item(18, 14).Value = -10368.75
$ws.Cells.Item(19, 8).Value = 2226
$ws.Cells.Item(19, 9).Value = 8
$ws.Cells.Item(19, 10).Value = 4444
$ws.Cells.Item(19, 11).Value = 8
$ws.Cells.Item(19, 12).Value = 4444
$ws.Cells.Item(19, 13).Value = 221
$ws.Cells.Item(19, 14).Value = -4902
$ws.Cells.Item(21, 8).Value = 6335
$ws.Cells.Item(21, 9).Value = 9448.333
$ws.Cells.Item(21, 10).Value = 2599
$ws.Cells.Item(21, 11).Value = 9448.333
$ws.Cells.Item(21, 12).Value = 2599
$ws.Cells.Item(21, 13).Value = -9074.333
$ws.Cells.Item(21, 14).Value = -3347
$ws.Cells.Item(22, 8).Value = 3335.3333
$ws.Cells.Item(22, 10).Value = 4995
$ws.Cells.Item(22, 12).Value = 4995
$ws.Cells.Item(22, 14).Value = -5593
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 13).ClearContents()
$ws.Cells.Item(25, 8).Value = 2351.625
$ws.Cells.Item(25, 9).Value = 1078.5
$ws.Cells.Item(25, 10).Value = 3624.75
$ws.Cells.Item(25, 11).Value = 1078.5
$ws.Cells.Item(25, 12).Value = 3624.75
$ws.Cells.Item(25, 13).Value = -676.5
$ws.Cells.Item(25, 14).Value = -4428.75
$ws.Cells.Item(27, 8).Value = 4495
$ws.Cells.Item(27, 10).Value = 4495
$ws.Cells.Item(27, 12).Value = 4495
$ws.Cells.Item(27, 14).Value = -4863
$ws.Cells.Item(30, 8).Value = 5998.3335
$ws.Cells.Item(30, 10).Value = 5998.3335
$ws.Cells.Item(30, 12).Value = 5998.3335
$ws.Cells.Item(30, 14).Value = -6298.3335
$ws.Cells.Item(32, 8).Value = 2354.6
$ws.Cells.Item(32, 9).Value = 2427.0212
$ws.Cells.Item(32, 10).Value = 1220
$ws.Cells.Item(32, 11).Value = 2427.0212
$ws.Cells.Item(32, 12).Value = 1220
$ws.Cells.Item(32, 13).Value = -2140.0212
$ws.Cells.Item(32, 14).Value = -1794
$ws.Cells.Item(35, 8).Value = 10388.5
$ws.Cells.Item(35, 9).Value = 9851.333
$ws.Cells.Item(35, 10).Value = 12000
$ws.Cells.Item(35, 11).Value = 9851.333
$ws.Cells.Item(35, 12).Value = 12000
$ws.Cells.Item(35, 13).Value = -9445.333
$ws.Cells.Item(35, 14).Value = -12812
$ws.Cells.Item(36, 8).Value = 9999.667
$ws.Cells.Item(36, 9).Value = 9999.5
$ws.Cells.Item(36, 10).Value = 10000
$ws.Cells.Item(36, 11).Value = 9999.5
$ws.Cells.Item(36, 12).Value = 10000
$ws.Cells.Item(36, 13).Value = -9653.5
$ws.Cells.Item(36, 14).Value = -10692
$ws.Cells.Item(37, 8).Value = 25555.111
$ws.Cells.Item(37, 10).Value = 25555.111
$ws.Cells.Item(37, 12).Value = 25555.111
$ws.Cells.Item(37, 14).Value = -26101.111
$ws.Cells.Item(39, 8).Value = 2395
$ws.Cells.Item(39, 9).Value = 2418.75
$ws.Cells.Item(39, 10).Value = 2300
$ws.Cells.Item(39, 11).Value = 2418.75
$ws.Cells.Item(39, 12).Value = 2300
$ws.Cells.Item(39, 13).Value = -1898.75
$ws.Cells.Item(39, 14).Value = -3340
$ws.Cells.Item(40, 8).Value = 34749.5
$ws.Cells.Item(40, 10).Value = 20000
$ws.Cells.Item(40, 12).Value = 20000
$ws.Cells.Item(40, 14).Value = -20352
$ws.Cells.Item(41, 8).Value = 12823
$ws.Cells.Item(41, 9).Value = 1108
$ws.Cells.Item(41, 10).Value = 32348
$ws.Cells.Item(41, 11).Value = 1108
$ws.Cells.Item(41, 12).Value = 32348
$ws.Cells.Item(41, 13).Value = -694
$ws.Cells.Item(41, 14).Value = -33176
$ws.Cells.Item(42, 8).Value = 42500
$ws.Cells.Item(42, 10).Value = 42500
$ws.Cells.Item(42, 12).Value = 42500
$ws.Cells.Item(42, 14).Value = -43472
$ws.Cells.Item(45, 8).Value = 8881.954
$ws.Cells.Item(45, 9).Value = 14971.2
$ws.Cells.Item(45, 11).Value = 14971.2
$ws.Cells.Item(45, 13).Value = -14594.2
$ws.Cells.Item(47, 8).Value = 30000
$ws.Cells.Item(47, 10).Value = 30000
$ws.Cells.Item(47, 12).Value = 30000
$ws.Cells.Item(47, 14).Value = -31450
$ws.Cells.Item(49, 8).Value = 25000
$ws.Cells.Item(49, 10).Value = 25000
$ws.Cells.Item(49, 12).Value = 25000
$ws.Cells.Item(49, 14).Value = -25520
$ws.Cells.Item(63, 8).Value = 3399.75
$ws.Cells.Item(63, 9).Value = 2943.889
$ws.Cells.Item(63, 11).Value = 2943.889
$ws.Cells.Item(63, 13).Value = -2257.889
$ws.Cells.Item(66, 8).Value = 3399.75
$ws.Cells.Item(66, 9).Value = 2943.889
$ws.Cells.Item(66, 11).Value = 14719.445
$ws.Cells.Item(66, 13).Value = -11287.445
$ws.Cells.Item(74, 8).Value = 9666.556
$ws.Cells.Item(74, 9).Value = 9868.692
$ws.Cells.Item(74, 11).Value = 9868.692
$ws.Cells.Item(74, 13).Value = -8994.692
$ws.Cells.Item(77, 8).Value = 9666.556
$ws.Cells.Item(77, 9).Value = 9868.692
$ws.Cells.Item(77, 11).Value = 49343.45999999999
$ws.Cells.Item(77, 13).Value = -44975.45999999999
$ws.Cells.Item(88, 8).Value = 1465.1875
$ws.Cells.Item(88, 10).Value = 1511
$ws.Cells.Item(88, 12).Value = 1511
$ws.Cells.Item(88, 14).Value = -2323
$ws.Cells.Item(91, 8).Value = 1465.1875
$ws.Cells.Item(91, 10).Value = 1511
$ws.Cells.Item(91, 12).Value = 1511
$ws.Cells.Item(91, 14).Value = -4319
$ws.Cells.Item(116, 8).Value = 698585.06
$ws.Cells.Item(116, 9).Value = 1021935.94
$ws.Cells.Item(116, 10).Value = 2137.077
$ws.Cells.Item(116, 11).Value = 1021935.94
$ws.Cells.Item(116, 12).Value = 2137.077
$ws.Cells.Item(116, 13).Value = -1019641.94
$ws.Cells.Item(116, 14).Value = -6725.077
$ws.Cells.Item(132, 8).Value = 2899.7805
$ws.Cells.Item(132, 9).Value = 2372.7778
$ws.Cells.Item(132, 10).Value = 6694.2
$ws.Cells.Item(132, 11).Value = 7118.3334
$ws.Cells.Item(132, 12).Value = 20082.6
$ws.Cells.Item(132, 13).Value = -4588.3334
$ws.Cells.Item(132, 14).Value = -25142.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 698585.06
$ws.Cells.Item(3, 9).Value = 1021935.94
$ws.Cells.Item(3, 10).Value = 2137.077
$ws.Cells.Item(3, 11).Value = 1021935.94
$ws.Cells.Item(3, 12).Value = 2137.077
$ws.Cells.Item(3, 13).Value = -1021821.94
$ws.Cells.Item(3, 14).Value = -2365.077
$ws.Cells.Item(22, 8).Value = 50337.5
$ws.Cells.Item(22, 10).Value = 125435.75
$ws.Cells.Item(22, 12).Value = 125435.75
$ws.Cells.Item(22, 14).Value = -125781.75
$ws.Cells.Item(55, 8).Value = 97123.5
$ws.Cells.Item(55, 10).Value = 97123.5
$ws.Cells.Item(55, 12).Value = 97123.5
$ws.Cells.Item(55, 14).Value = -97669.5
$ws.Cells.Item(94, 8).Value = 2952.125
$ws.Cells.Item(94, 9).Value = 2945.8572
$ws.Cells.Item(94, 11).Value = 2945.8572
$ws.Cells.Item(94, 13).Value = -2494.8572
$ws.Cells.Item(105, 8).Value = 4831.8184
$ws.Cells.Item(105, 9).Value = 3234.7058
$ws.Cells.Item(105, 11).Value = 3234.7058
$ws.Cells.Item(105, 13).Value = -1487.7058
$ws.Cells.Item(107, 8).Value = 1854.16
$ws.Cells.Item(107, 9).Value = 1868.4667
$ws.Cells.Item(107, 11).Value = 1868.4667
$ws.Cells.Item(107, 13).Value = 51.53330000000005

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 111114664
$ws.Cells.Item(62, 9).Value = 333335650
$ws.Cells.Item(62, 10).Value = 4174.6665
$ws.Cells.Item(62, 11).Value = 333335650
$ws.Cells.Item(62, 12).Value = 4174.6665
$ws.Cells.Item(62, 13).Value = -333335026
$ws.Cells.Item(62, 14).Value = -5422.6665
$ws.Cells.Item(65, 8).Value = 111114664
$ws.Cells.Item(65, 9).Value = 333335650
$ws.Cells.Item(65, 10).Value = 4174.6665
$ws.Cells.Item(65, 11).Value = 1666678250
$ws.Cells.Item(65, 12).Value = 20873.3325
$ws.Cells.Item(65, 13).Value = -1666675130
$ws.Cells.Item(65, 14).Value = -27113.3325
$ws.Cells.Item(93, 8).Value = 3855.4285
$ws.Cells.Item(93, 9).Value = 3855.4285
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 3855.4285
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = -1983.4285
$ws.Cells.Item(93, 14).ClearContents()
$ws.Cells.Item(121, 8).Value = 64500
$ws.Cells.Item(121, 10).Value = 64500
$ws.Cells.Item(121, 12).Value = 64500
$ws.Cells.Item(121, 14).Value = -67120
$ws.Cells.Item(125, 8).Value = 59296
$ws.Cells.Item(125, 10).Value = 59296
$ws.Cells.Item(125, 12).Value = 59296
$ws.Cells.Item(125, 14).Value = -64216
$ws.Cells.Item(132, 8).Value = 5308.125
$ws.Cells.Item(132, 10).Value = 21268.166
$ws.Cells.Item(132, 12).Value = 63804.49800000001
$ws.Cells.Item(132, 14).Value = -68864.498
$ws.Cells.Item(138, 8).Value = 180005.75
$ws.Cells.Item(138, 10).Value = 180005.75
$ws.Cells.Item(138, 12).Value = 180005.75
$ws.Cells.Item(138, 14).Value = -190285.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(25, 8).Value = 1150
$ws.Cells.Item(25, 9).Value = 1150
$ws.Cells.Item(25, 11).Value = 3450
$ws.Cells.Item(25, 13).Value = -3281
$ws.Cells.Item(30, 8).Value = 1150
$ws.Cells.Item(30, 9).Value = 1150
$ws.Cells.Item(30, 11).Value = 3450
$ws.Cells.Item(30, 13).Value = -3348
$ws.Cells.Item(64, 8).Value = 2750.625
$ws.Cells.Item(64, 9).Value = 1613.5
$ws.Cells.Item(64, 10).Value = 3129.6667
$ws.Cells.Item(64, 11).Value = 4840.5
$ws.Cells.Item(64, 12).Value = 9389.000100000001
$ws.Cells.Item(64, 13).Value = -4570.5
$ws.Cells.Item(64, 14).Value = -9929.000100000001
$ws.Cells.Item(67, 8).Value = 2750.625
$ws.Cells.Item(67, 9).Value = 1613.5
$ws.Cells.Item(67, 10).Value = 3129.6667
$ws.Cells.Item(67, 11).Value = 4840.5
$ws.Cells.Item(67, 12).Value = 9389.000100000001
$ws.Cells.Item(67, 13).Value = -3904.5
$ws.Cells.Item(67, 14).Value = -11261.0001
$ws.Cells.Item(68, 8).Value = 3332.8333
$ws.Cells.Item(68, 10).Value = 3332.8333
$ws.Cells.Item(68, 12).Value = 9998.499899999999
$ws.Cells.Item(68, 14).Value = -11620.4999
$ws.Cells.Item(71, 8).Value = 3332.8333
$ws.Cells.Item(71, 10).Value = 3332.8333
$ws.Cells.Item(71, 12).Value = 29995.4997
$ws.Cells.Item(71, 14).Value = -38107.4997
$ws.Cells.Item(107, 8).Value = 527.5714
$ws.Cells.Item(107, 9).Value = 504.75
$ws.Cells.Item(107, 11).Value = 1514.25
$ws.Cells.Item(107, 13).Value = 405.75
$ws.Cells.Item(125, 8).Value = 11682

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 40935710
$ws.Cells.Item(2, 9).Value = 55555584
$ws.Cells.Item(2, 11).Value = 55555584
$ws.Cells.Item(2, 13).Value = -55555471
$ws.Cells.Item(70, 8).Value = 8662
$ws.Cells.Item(70, 9).Value = 9900
$ws.Cells.Item(70, 10).Value = 8043
$ws.Cells.Item(70, 11).Value = 9900
$ws.Cells.Item(70, 12).Value = 8043
$ws.Cells.Item(70, 13).Value = -9630
$ws.Cells.Item(70, 14).Value = -8583
$ws.Cells.Item(73, 8).Value = 8662
$ws.Cells.Item(73, 9).Value = 9900
$ws.Cells.Item(73, 10).Value = 8043
$ws.Cells.Item(73, 11).Value = 9900
$ws.Cells.Item(73, 12).Value = 8043
$ws.Cells.Item(73, 13).Value = -8964
$ws.Cells.Item(73, 14).Value = -9915
$ws.Cells.Item(80, 8).Value = 17695114
$ws.Cells.Item(80, 9).Value = 21907508
$ws.Cells.Item(80, 10).Value = 3066.4
$ws.Cells.Item(80, 11).Value = 21907508
$ws.Cells.Item(80, 12).Value = 3066.4
$ws.Cells.Item(80, 13).Value = -21906510
$ws.Cells.Item(80, 14).Value = -5062.4
$ws.Cells.Item(83, 8).Value = 17695114
$ws.Cells.Item(83, 9).Value = 21907508
$ws.Cells.Item(83, 10).Value = 3066.4
$ws.Cells.Item(83, 11).Value = 109537540
$ws.Cells.Item(83, 12).Value = 15332
$ws.Cells.Item(83, 13).Value = -109532548
$ws.Cells.Item(83, 14).Value = -25316
$ws.Cells.Item(102, 8).Value = 4119.207
$ws.Cells.Item(102, 9).Value = 4697.737
$ws.Cells.Item(102, 10).Value = 3020
$ws.Cells.Item(102, 11).Value = 4697.737
$ws.Cells.Item(102, 12).Value = 3020
$ws.Cells.Item(102, 13).Value = -3075.737
$ws.Cells.Item(102, 14).Value = -6264
$ws.Cells.Item(110, 8).Value = 145849.5
$ws.Cells.Item(110, 10).Value = 145849.5
$ws.Cells.Item(110, 12).Value = 145849.5
$ws.Cells.Item(110, 14).Value = -154029.5
$ws.Cells.Item(122, 8).Value = 6756.9165
$ws.Cells.Item(122, 9).Value = 6363.45
$ws.Cells.Item(122, 11).Value = 19090.35
$ws.Cells.Item(122, 13).Value = -16640.35
$ws.Cells.Item(140, 8).Value = 52638.777
$ws.Cells.Item(140, 10).Value = 57343.625
$ws.Cells.Item(140, 12).Value = 57343.625
$ws.Cells.Item(140, 14).Value = -67703.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 10408.637
$ws.Cells.Item(7, 10).Value = 19333.334
$ws.Cells.Item(7, 12).Value = 19333.334
$ws.Cells.Item(7, 14).Value = -19557.334
$ws.Cells.Item(22, 8).Value = 755.5
$ws.Cells.Item(22, 9).Value = 755.5
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 755.5
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -460.5
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(27, 8).Value = 755.5
$ws.Cells.Item(27, 9).Value = 755.5
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 755.5
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = -648.5
$ws.Cells.Item(27, 14).ClearContents()
$ws.Cells.Item(46, 8).Value = 3427.25
$ws.Cells.Item(46, 9).Value = 1828.1538
$ws.Cells.Item(46, 11).Value = 1828.1538
$ws.Cells.Item(46, 13).Value = -1640.1538
$ws.Cells.Item(61, 8).Value = 9595.529
$ws.Cells.Item(61, 9).Value = 11283.923
$ws.Cells.Item(61, 11).Value = 11283.923
$ws.Cells.Item(61, 13).Value = -11081.923
$ws.Cells.Item(97, 8).Value = 45179
$ws.Cells.Item(97, 10).Value = 45179
$ws.Cells.Item(97, 12).Value = 45179
$ws.Cells.Item(97, 14).Value = -47161
$ws.Cells.Item(110, 8).Value = 45000
$ws.Cells.Item(110, 10).Value = 45000
$ws.Cells.Item(110, 12).Value = 45000
$ws.Cells.Item(110, 14).Value = -53180
$ws.Cells.Item(113, 8).Value = 9595.529
$ws.Cells.Item(113, 9).Value = 11283.923
$ws.Cells.Item(113, 11).Value = 11283.923
$ws.Cells.Item(113, 13).Value = -9113.923
$ws.Cells.Item(122, 8).Value = 4584.857
$ws.Cells.Item(122, 9).Value = 4850
$ws.Cells.Item(122, 10).Value = 4231.3335
$ws.Cells.Item(122, 11).Value = 14550
$ws.Cells.Item(122, 12).Value = 12694.0005
$ws.Cells.Item(122, 13).Value = -12100
$ws.Cells.Item(122, 14).Value = -17594.0005
$ws.Cells.Item(126, 8).Value = 10408.637
$ws.Cells.Item(126, 10).Value = 19333.334
$ws.Cells.Item(126, 12).Value = 58000.00199999999
$ws.Cells.Item(126, 14).Value = -62940.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 113588.11
$ws.Cells.Item(96, 9).Value = 169515.67
$ws.Cells.Item(96, 10).Value = 1733
$ws.Cells.Item(96, 11).Value = 169515.67
$ws.Cells.Item(96, 12).Value = 1733
$ws.Cells.Item(96, 13).Value = -168142.67
$ws.Cells.Item(96, 14).Value = -4479
$ws.Cells.Item(132, 8).Value = 3505.5667
$ws.Cells.Item(132, 9).Value = 2791.6785
$ws.Cells.Item(132, 11).Value = 8375.0355
$ws.Cells.Item(132, 13).Value = -5845.0355
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
